# Adds a new literature-review row (row 14) to the "data" worksheet,
# recording the Depicker et al., 2021 source, and moves the active
# selection to the next empty cell (E14), matching the author's workflow
# of having just finished filling in the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("A14").Value = "Depicker et al., 2021"
$ws.Range("B14").Value = "Kivu Rift, Africa (DRC, Rwanda and Burundi)"
$ws.Range("C14").Value = "deforestation"
$ws.Range("D14").Value = "Modelling"

$ws.Activate() | Out-Null
$ws.Range("E14").Select() | Out-Null
